$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $r = $ws.Range($cellAddr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "44.095.67"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "2.255.53"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "258.38"
$ws.Range("E5").Value = "  +3.07%  "
Set-TextValue "D6" "79.46"
$ws.Range("E6").Value = "  +6.53%  "
Set-TextValue "D7" "0.627"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E8").Value = "  -0.16%  "
Set-TextValue "D9" "0.609"
$ws.Range("E9").Value = "  +3.11%  "
Set-TextValue "D10" "43.50"
$ws.Range("E10").Value = "  +7.61%  "
Set-TextValue "D11" "0.0930"
$ws.Range("E11").Value = "  +1.28%  "
Set-TextValue "D12" "7.15"
$ws.Range("E12").Value = "  +4.63%  "
Set-TextValue "D13" "0.104"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "2.575.99"
$ws.Range("E14").Value = "  +1.26%  "
Set-TextValue "D15" "14.80"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "2.277.32"
$ws.Range("E16").Value = "  +3.12%  "
Set-TextValue "D17" "0.798"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "43.949.68"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("E19").Value = "  +2.35%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D20" "6.10"
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D21" "71.71"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  +7.32%  "
Set-TextValue "D23" "234.99"
$ws.Range("E23").Value = "  +2.57%  "
Set-TextValue "D24" "9.55"
$ws.Range("E24").Value = "  +1.40%  "
Set-TextValue "D26" "42.87"
$ws.Range("E26").Value = "  +9.54%  "
Set-TextValue "D27" "10.86"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  -0.91%  "
Set-TextValue "D29" "2.24"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  -1.04%  "
Set-TextValue "D31" "173.77"
$ws.Range("E31").Value = "  +2.27%  "
Set-TextValue "D32" "20.77"
$ws.Range("E32").Value = "  +3.15%  "
Set-TextValue "D33" "0.0882"
$ws.Range("E33").Value = "  +11.10%  "
Set-TextValue "D34" "5.34"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E35").Value = "  +2.09%  "
Set-TextValue "D36" "0.0367"
$ws.Range("E36").Value = "  +12.39%  "
Set-TextValue "D37" "4.55"
$ws.Range("E37").Value = "  +3.50%  "
Set-TextValue "D38" "0.110"
$ws.Range("E38").Value = "  +0.22%  "
Set-TextValue "D39" "13.38"
$ws.Range("E39").Value = "  +11.37%  "
$ws.Range("E40").Value = "  +20.35%  "
Set-TextValue "D41" "2.15"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.205"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D43" "62.39"
$ws.Range("E43").Value = "  +5.94%  "
Set-TextValue "D44" "5.42"
$ws.Range("E44").Value = "  +2.82%  "
Set-TextValue "D45" "105.70"
$ws.Range("E45").Value = "  +3.18%  "
Set-TextValue "D46" "8.59"
$ws.Range("E46").Value = "  +0.85%  "
Set-TextValue "D47" "0.478"
$ws.Range("E47").Value = "  +0.40%  "
Set-TextValue "D48" "0.0991"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("E51").Value = "  +26.98%  "
